$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.414.46"
$ws.Range("E2").Value = "  +5.35%  "
$ws.Range("D3").Value = "1.809.93"
$ws.Range("E3").Value = "  +4.22%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'316.56"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'0.5540"
$ws.Range("E7").Value = "  +10.81%  "
$ws.Range("D8").Value = "'0.3854"
$ws.Range("E8").Value = "  +9.54%  "
$ws.Range("D9").Value = "'0.07592"
$ws.Range("E9").Value = "  +4.73%  "
$ws.Range("D10").Value = "'42.94"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").Value = "'1.132"
$ws.Range("E11").Value = "  +7.08%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "'21.16"
$ws.Range("E13").Value = "  +5.58%  "
$ws.Range("D14").Value = "'6.239"
$ws.Range("E14").Value = "  +5.00%  "
$ws.Range("D15").Value = "1.810.11"
$ws.Range("E15").Value = "  +4.48%  "
$ws.Range("D16").Value = "'7.313"
$ws.Range("E16").Value = "  +6.43%  "
$ws.Range("D17").Value = "'91.30"
$ws.Range("E17").Value = "  +5.76%  "
$ws.Range("D18").Value = "'0.00001073"
$ws.Range("E18").Value = "  +3.85%  "
$ws.Range("D19").Value = "'0.06479"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "'17.25"
$ws.Range("E21").Value = "  +3.82%  "
$ws.Range("D22").Value = "'5.995"
$ws.Range("E22").Value = "  +4.57%  "
$ws.Range("D23").Value = "28.425.31"
$ws.Range("E23").Value = "  +5.04%  "
$ws.Range("D24").Value = "'11.31"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("D25").Value = "'2.124"
$ws.Range("E25").Value = "  +3.34%  "
$ws.Range("D26").Value = "'157.43"
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("E27").Value = "  +3.82%  "
$ws.Range("D28").Value = "'2.438"
$ws.Range("E28").Value = "  +15.61%  "
$ws.Range("D29").Value = "2.017.84"
$ws.Range("E29").Value = "  +4.25%  "
$ws.Range("D30").Value = "'124.02"
$ws.Range("E30").Value = "  +2.94%  "
$ws.Range("D31").Value = "'1.164"
$ws.Range("E31").Value = "  +9.43%  "
$ws.Range("D32").Value = "'0.1039"
$ws.Range("E32").Value = "  +9.33%  "
$ws.Range("D33").Value = "'5.768"
$ws.Range("E33").Value = "  +7.11%  "
$ws.Range("D34").Value = "'3.644"
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("D35").Value = "'0.2275"
$ws.Range("E35").Value = "  +14.16%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'8.900"
$ws.Range("E36").Value = "  +19.08%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02324"
$ws.Range("E37").Value = "  +6.08%  "
$ws.Range("D38").Value = "'0.06197"
$ws.Range("E38").Value = "  +4.41%  "
$ws.Range("D39").Value = "'11.63"
$ws.Range("E39").Value = "  +5.48%  "
$ws.Range("D40").Value = "'0.6391"
$ws.Range("E40").Value = "  +6.46%  "
$ws.Range("D41").Value = "'5.027"
$ws.Range("E41").Value = "  +5.61%  "
$ws.Range("D42").Value = "'1.182"
$ws.Range("E42").Value = "  +6.40%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  -3.25%  "
$ws.Range("D45").Value = "'13.47"
$ws.Range("E45").Value = "  +5.01%  "
$ws.Range("D46").Value = "'0.5998"
$ws.Range("E46").Value = "  +6.40%  "
$ws.Range("D47").Value = "'3.696"
$ws.Range("E47").Value = "  +3.36%  "
$ws.Range("D48").Value = "'123.23"
$ws.Range("E48").Value = "  +3.15%  "
$ws.Range("D49").Value = "'1.973"
$ws.Range("E49").Value = "  +6.60%  "
$ws.Range("D50").Value = "'1.144"
$ws.Range("E50").Value = "  +4.28%  "
$ws.Range("D51").Value = "'0.06928"
$ws.Range("E51").Value = "  +4.04%  "
